# Updates the cryptos list with refreshed price / 1h-volume data
# (swap of the Binance-PegBSC-USD / InternetComputer(DFINITY) rows included)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.189.00'
$ws.Range('E2').Value = '  +0.40%  '

$ws.Range('D3').Value = '2.661.41'
$ws.Range('E3').Value = '  +3.18%  '

$ws.Range('E4').Value = '  +0.10%  '

$ws.Range('D5').Value = '''607.32'
$ws.Range('E5').Value = '  +4.29%  '

$ws.Range('D6').Value = '''142.92'
$ws.Range('E6').Value = '  -1.07%  '

$ws.Range('E7').Value = '  +0.12%  '

$ws.Range('E8').Value = '  -1.21%  '

$ws.Range('D9').Value = '2.661.56'
$ws.Range('E9').Value = '  +3.20%  '

$ws.Range('E10').Value = '  -0.12%  '

$ws.Range('E11').Value = '  +1.47%  '

$ws.Range('E12').Value = '  +0.72%  '

$ws.Range('D13').Value = '''0.356'
$ws.Range('E13').Value = '  +1.90%  '

$ws.Range('D14').Value = '''27.27'
$ws.Range('E14').Value = '  +0.87%  '

$ws.Range('D15').Value = '3.143.28'
$ws.Range('E15').Value = '  +3.35%  '

$ws.Range('D16').Value = '63.061.85'
$ws.Range('E16').Value = '  +0.38%  '

$ws.Range('E17').Value = '  -0.39%  '

$ws.Range('D18').Value = '2.637.85'
$ws.Range('E18').Value = '  +2.30%  '

$ws.Range('D19').Value = '''11.44'
$ws.Range('E19').Value = '  +3.11%  '

$ws.Range('D20').Value = '''339.31'
$ws.Range('E20').Value = '  -0.41%  '

$ws.Range('E21').Value = '  +1.05%  '

$ws.Range('D22').Value = '''6.84'
$ws.Range('E22').Value = '  +3.05%  '

$ws.Range('D23').Value = '''1.00'
$ws.Range('E23').Value = '  +0.12%  '

$ws.Range('D24').Value = '''67.60'
$ws.Range('E24').Value = '  +0.34%  '

$ws.Range('E25').Value = '  +2.84%  '

$ws.Range('D26').Value = '''1.56'
$ws.Range('E26').Value = '  -2.16%  '

$ws.Range('E27').Value = '  -0.36%  '

$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '''8.47'
$ws.Range('E28').Value = '  +2.62%  '

$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  -0.07%  '

$ws.Range('D30').Value = '''537.56'
$ws.Range('E30').Value = '  +16.69%  '

$ws.Range('D31').Value = '''7.84'
$ws.Range('E31').Value = '  -2.13%  '

$ws.Range('D32').Value = '''2.03'
$ws.Range('E32').Value = '  +5.34%  '

$ws.Range('E33').Value = '  +9.14%  '

$ws.Range('D34').Value = '0.0₃0807'
$ws.Range('E34').Value = '  +1.07%  '

$ws.Range('D35').Value = '''172.72'
$ws.Range('E35').Value = '  -1.96%  '

$ws.Range('D36').Value = '''5.08'
$ws.Range('E36').Value = '  +13.65%  '

$ws.Range('E37').Value = '  -0.01%  '

$ws.Range('E38').Value = '  +2.10%  '

$ws.Range('E39').Value = '  +1.83%  '

$ws.Range('E40').Value = '  +7.36%  '

$ws.Range('D41').Value = '''174.39'
$ws.Range('E41').Value = '  +10.21%  '

$ws.Range('E42').Value = '  +0.04%  '

$ws.Range('E43').Value = '  +1.33%  '

$ws.Range('D44').Value = '''21.99'
$ws.Range('E44').Value = '  +3.96%  '

$ws.Range('E45').Value = '  +4.49%  '

$ws.Range('E46').Value = '  -0.73%  '

$ws.Range('E47').Value = '  +0.11%  '

$ws.Range('E48').Value = '  +1.30%  '

$ws.Range('D49').Value = '''18.71'
$ws.Range('E49').Value = '  +4.12%  '

$ws.Range('D51').Value = '''11.34'
$ws.Range('E51').Value = '  -0.59%  '
